# CardiologyData.xlsx — fill in Row 3 (% missing), Row 7 (Mode) and Row 9 (Standard Deviation)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3: % of missing values in the data ----
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.0064
$ws.Range("F3").NumberFormat = "0.00%"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.0097
$ws.Range("H3").NumberFormat = "0.00%"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0.0064
$ws.Range("P3").NumberFormat = "0.00%"

# ---- Row 7: Mode ----
$ws.Range("B7").Value = 58
$ws.Range("C7").Value = "Male"
$ws.Range("D7").Value = "NoTang"
$ws.Range("E7").Value = 120
$ws.Range("F7").Value = 204
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = "Normal"
$ws.Range("I7").Value = 76.8
$ws.Range("J7").Value = 162
$ws.Range("K7").Value = $false
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = "Up"
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = "Normal"
$ws.Range("P7").Value = "Healthy"

# ---- Row 9: Standard Deviation ----
$ws.Range("B9").Value = 9.775942
$ws.Range("C9").Value = "NA"
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = 17.49002
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = 0.3537851
$ws.Range("H9").Value = "NA"
$ws.Range("I9").Value = 10.71816
$ws.Range("J9").Value = 22.76805
$ws.Range("K9").Value = 0.4690154
$ws.Range("L9").Value = 1.166037
$ws.Range("M9").Value = "NA"
$ws.Range("N9").Value = 0.9312189
$ws.Range("O9").Value = "NA"
$ws.Range("P9").Value = "NA"

# Return selection to the sheet's default cell
$ws.Range("A1").Select() | Out-Null
